$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cxcl12"
$ws.Cells.Item(2, 3).Value = "Cd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 130.955829
$ws.Cells.Item(2, 8).Value = 392.867487
$ws.Cells.Item(2, 9).Value = 0.5336535908353144
$ws.Cells.Item(2, 10).Value = 0.5336535908353144
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 1.635729666666667
$ws.Cells.Item(2, 14).Value = 4.907189
$ws.Cells.Item(2, 15).Value = 0.5314629201652572
$ws.Cells.Item(2, 16).Value = 0.5314629201652572
$ws.Cells.Item(2, 17).Value = 214.208334518227
$ws.Cells.Item(2, 18).Value = 1927.875010664043
$ws.Cells.Item(2, 19).Value = 0.2836170957420115
$ws.Cells.Item(2, 20).Value = 0.2836170957420115

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cxcl12"
$ws.Cells.Item(3, 3).Value = "Cd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 130.955829
$ws.Cells.Item(3, 8).Value = 392.867487
$ws.Cells.Item(3, 9).Value = 0.5336535908353144
$ws.Cells.Item(3, 10).Value = 0.5336535908353144
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.170126666666667
$ws.Cells.Item(3, 14).Value = 3.51038
$ws.Cells.Item(3, 15).Value = 0.3801844203860328
$ws.Cells.Item(3, 16).Value = 0.3801844203860328
$ws.Cells.Item(3, 17).Value = 153.23490766834
$ws.Cells.Item(3, 18).Value = 1379.11416901506
$ws.Cells.Item(3, 19).Value = 0.2028867811186491
$ws.Cells.Item(3, 20).Value = 0.2028867811186491

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cxcl12"
$ws.Cells.Item(4, 3).Value = "Cd4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 130.955829
$ws.Cells.Item(4, 8).Value = 392.867487
$ws.Cells.Item(4, 9).Value = 0.5336535908353144
$ws.Cells.Item(4, 10).Value = 0.5336535908353144
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.2719306666666667
$ws.Cells.Item(4, 14).Value = 0.8157920000000001
$ws.Cells.Item(4, 15).Value = 0.08835265944870997
$ws.Cells.Item(4, 16).Value = 0.08835265944870996
$ws.Cells.Item(4, 17).Value = 35.61090588385601
$ws.Cells.Item(4, 18).Value = 320.498152954704
$ws.Cells.Item(4, 19).Value = 0.04714971397465374
$ws.Cells.Item(4, 20).Value = 0.04714971397465374

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cxcl12"
$ws.Cells.Item(5, 3).Value = "Cd4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 66.39541
$ws.Cells.Item(5, 8).Value = 199.18623
$ws.Cells.Item(5, 9).Value = 0.2705656497465488
$ws.Cells.Item(5, 10).Value = 0.2705656497465488
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 1.635729666666667
$ws.Cells.Item(5, 14).Value = 4.907189
$ws.Cells.Item(5, 15).Value = 0.5314629201652572
$ws.Cells.Item(5, 16).Value = 0.5314629201652572
$ws.Cells.Item(5, 17).Value = 108.6049418674966
$ws.Cells.Item(5, 18).Value = 977.44447680747
$ws.Cells.Item(5, 19).Value = 0.143795610310711
$ws.Cells.Item(5, 20).Value = 0.143795610310711

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cxcl12"
$ws.Cells.Item(6, 3).Value = "Cd4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 66.39541
$ws.Cells.Item(6, 8).Value = 199.18623
$ws.Cells.Item(6, 9).Value = 0.2705656497465488
$ws.Cells.Item(6, 10).Value = 0.2705656497465488
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.170126666666667
$ws.Cells.Item(6, 14).Value = 3.51038
$ws.Cells.Item(6, 15).Value = 0.3801844203860328
$ws.Cells.Item(6, 16).Value = 0.3801844203860328
$ws.Cells.Item(6, 17).Value = 77.69103978526667
$ws.Cells.Item(6, 18).Value = 699.2193580674
$ws.Cells.Item(6, 19).Value = 0.102864844725262
$ws.Cells.Item(6, 20).Value = 0.102864844725262

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cxcl12"
$ws.Cells.Item(7, 3).Value = "Cd4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 66.39541
$ws.Cells.Item(7, 8).Value = 199.18623
$ws.Cells.Item(7, 9).Value = 0.2705656497465488
$ws.Cells.Item(7, 10).Value = 0.2705656497465488
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.2719306666666667
$ws.Cells.Item(7, 14).Value = 0.8157920000000001
$ws.Cells.Item(7, 15).Value = 0.08835265944870997
$ws.Cells.Item(7, 16).Value = 0.08835265944870996
$ws.Cells.Item(7, 17).Value = 18.05494810490667
$ws.Cells.Item(7, 18).Value = 162.49453294416
$ws.Cells.Item(7, 19).Value = 0.02390519471057577
$ws.Cells.Item(7, 20).Value = 0.02390519471057577

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Cxcl12"
$ws.Cells.Item(8, 3).Value = "Cd4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 48.043585
$ws.Cells.Item(8, 8).Value = 144.130755
$ws.Cells.Item(8, 9).Value = 0.1957807594181367
$ws.Cells.Item(8, 10).Value = 0.1957807594181367
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 1.635729666666667
$ws.Cells.Item(8, 14).Value = 4.907189
$ws.Cells.Item(8, 15).Value = 0.5314629201652572
$ws.Cells.Item(8, 16).Value = 0.5314629201652572
$ws.Cells.Item(8, 17).Value = 78.58631727752166
$ws.Cells.Item(8, 18).Value = 707.2768554976949
$ws.Cells.Item(8, 19).Value = 0.1040502141125346
$ws.Cells.Item(8, 20).Value = 0.1040502141125346

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Cxcl12"
$ws.Cells.Item(9, 3).Value = "Cd4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 48.043585
$ws.Cells.Item(9, 8).Value = 144.130755
$ws.Cells.Item(9, 9).Value = 0.1957807594181367
$ws.Cells.Item(9, 10).Value = 0.1957807594181367
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.170126666666667
$ws.Cells.Item(9, 14).Value = 3.51038
$ws.Cells.Item(9, 15).Value = 0.3801844203860328
$ws.Cells.Item(9, 16).Value = 0.3801844203860328
$ws.Cells.Item(9, 17).Value = 56.21707997076667
$ws.Cells.Item(9, 18).Value = 505.9537197369
$ws.Cells.Item(9, 19).Value = 0.07443279454212165
$ws.Cells.Item(9, 20).Value = 0.07443279454212162

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Cxcl12"
$ws.Cells.Item(10, 3).Value = "Cd4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 48.043585
$ws.Cells.Item(10, 8).Value = 144.130755
$ws.Cells.Item(10, 9).Value = 0.1957807594181367
$ws.Cells.Item(10, 10).Value = 0.1957807594181367
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.2719306666666667
$ws.Cells.Item(10, 14).Value = 0.8157920000000001
$ws.Cells.Item(10, 15).Value = 0.08835265944870997
$ws.Cells.Item(10, 16).Value = 0.08835265944870996
$ws.Cells.Item(10, 17).Value = 13.06452409810667
$ws.Cells.Item(10, 18).Value = 117.58071688296
$ws.Cells.Item(10, 19).Value = 0.01729775076348045
$ws.Cells.Item(10, 20).Value = 0.01729775076348045
